# Update results and plots
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column J: "Binary Bird Dataset" results, mirroring the layout of column H ("DESED 2022")
$ws.Range("J1").Value = "Binary Bird Dataset"
$ws.Range("J1").Style = "Normal"
$ws.Cells.Item(1, 10).Font.Bold = $true

$ws.Range("J3").Value = "0.4154"
$ws.Range("J4").Value = "0.4154"
$ws.Range("J5").Value = "1.0"
$ws.Range("J6").Value = "0.5869"
$ws.Cells.Item(6, 10).Font.Bold = $true

$ws.Range("J8").Value = "0.8834"
$ws.Range("J9").Value = "0.9163"
$ws.Range("J10").Value = "0.7915"
$ws.Range("J11").Value = "0.8493"
$ws.Cells.Item(11, 10).Font.Bold = $true

# Column width for the new column
$ws.Columns.Item(10).ColumnWidth = 18

# H11 loses its (redundant) explicit style, becoming default/plain like the rest of that column
$ws.Cells.Item(11, 8).Font.Bold = $false

# Move the active selection to the new last-used cell, like the author did
$ws.Range("J6").Select()
